$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.052.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.122.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.66"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.366"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.50%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.134.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.729"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.202"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.69"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.271.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.732.61"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.137.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +8.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.77"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000200"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "435.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.66"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.53"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.339.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +44.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.227"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +18.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.20"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.168"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.61"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.88"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.43%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.92"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.48%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "494.21"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.32"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.442"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +10.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.40"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -8.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.08"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.43"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.19%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.701"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.35"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.50%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.36"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.29%  "
